$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("展览")
$ws.Range("B2").NumberFormat = "@"
$ws.Range("B2").Value = "2024-07-27"
$ws.Range("C2").Value = "广州·AP动漫游戏嘉年华"
$ws.Range("D2").Value = "新港东路630-638号 南丰国际会展中心"
$ws.Range("E2").Value = "2024.07.27 09:00-07.28 17:00"
$ws.Range("F2").Value = 10370
$ws.Range("G2").Value = 80
$ws.Range("H2").Value = "https://show.bilibili.com/platform/detail.html?id=87213"
$ws.Range("I2").Value = "//i1.hdslb.com/bfs/openplatform/202406/3Z8rGZPP1718164976101.jpeg"
$ws.Range("C3").Value = "广州·AP动漫游戏嘉年华·徐慧内场"
$ws.Range("E3").Value = "2024.07.27 09:00-07.27 17:00"
$ws.Range("F3").Value = 242
$ws.Range("G3").Value = 158
$ws.Range("H3").Value = "https://show.bilibili.com/platform/detail.html?id=87801"
$ws.Range("I3").Value = "//i2.hdslb.com/bfs/openplatform/202406/He5WJVKl1718701764045.jpeg"
$ws.Range("C4").Value = "广州·LookLook动漫嘉年华2th"
$ws.Range("D4").Value = "展贸东路200号 恒达智慧汽车城"
$ws.Range("E4").Value = "2024.07.27 10:00-07.28 17:30"
$ws.Range("F4").Value = 1141
$ws.Range("G4").Value = 68
$ws.Range("H4").Value = "https://show.bilibili.com/platform/detail.html?id=87217"
$ws.Range("I4").Value = "//i2.hdslb.com/bfs/openplatform/202407/wjvAqamr1720170199991.jpeg"
$ws.Range("C5").Value = "广州·原神x星穹x崩only"
$ws.Range("D5").Value = "鸿盛二路巨大创意产业园 巨大产业园·智汇港"
$ws.Range("E5").Value = "2024.07.27 10:00-07.27 17:00"
$ws.Range("F5").Value = 1027
$ws.Range("G5").Value = 55
$ws.Range("H5").Value = "https://show.bilibili.com/platform/detail.html?id=87184"
$ws.Range("I5").Value = "//i1.hdslb.com/bfs/openplatform/202406/u67hjpFi1718160712051.jpeg"
$ws.Range("C6").Value = "广州·咒术回战only 1th"
$ws.Range("D6").Value = "奥体南路12号 优托邦(奥体旗舰店)"
$ws.Range("F6").Value = 810
$ws.Range("G6").Value = 45
$ws.Range("H6").Value = "https://show.bilibili.com/platform/detail.html?id=87842"
$ws.Range("I6").Value = "//i0.hdslb.com/bfs/openplatform/202406/uw0jbGu21718943362522.jpeg"
$ws.Range("C7").Value = "广州·畅响未来MIKU PARTY 同人&谷子市集 【免费入场】"
$ws.Range("D7").Value = "林和中路63号 东方宝泰购物广场"
$ws.Range("E7").Value = "2024.07.27 12:00-07.28 21:00"
$ws.Range("F7").Value = 255
$ws.Range("G7").Value = 20
$ws.Range("H7").Value = "https://show.bilibili.com/platform/detail.html?id=88767"
$ws.Range("I7").Value = "//i1.hdslb.com/bfs/openplatform/202407/HmrOWl5F1720176007755.jpeg"
$ws.Range("C8").Value = "广州·蔚蓝档案only（取消）"
$ws.Range("D8").Value = "奥体南路12号 优托邦(奥体旗舰店)"
$ws.Range("E8").Value = "2024.07.27 10:00-07.27 17:00"
$ws.Range("F8").Value = 58
$ws.Range("G8").Value = "不可售"
$ws.Range("H8").Value = "https://show.bilibili.com/platform/detail.html?id=87534"
$ws.Range("I8").Value = "//i1.hdslb.com/bfs/openplatform/202406/EBpwLeYi1718355256664.jpeg"
$ws.Range("F9").Value = 329
$ws.Range("F10").Value = 1078
$ws.Range("F12").Value = 143
$ws.Range("F13").Value = 833
$ws.Range("F14").Value = 362
$ws.Range("F15").Value = 1742
$ws.Range("F17").Value = 877
$ws.Range("F18").Value = 787
$ws.Range("F19").Value = 520
$ws.Range("F20").Value = 757
$ws.Range("F21").Value = 854
$ws.Range("F24").Value = 83
$ws.Range("F25").Value = 591
$ws.Range("F26").Value = 598
$ws.Range("F29").Value = 988
$ws.Range("F31").Value = 468
$ws.Range("F32").Value = 137
$ws.Range("F33").Value = 238
$ws.Range("F34").Value = 205
$ws.Range("F35").Value = 524
$ws.Range("F36").Value = 1569
$ws.Range("F37").Value = 347
$ws.Range("F39").Value = 1358
$ws.Range("F40").Value = 395
$ws.Range("F41").Value = 112
$ws.Range("F42").Value = 41
$ws.Range("F43").Value = 69
$ws.Range("F45").Value = 61
$ws.Range("F46").Value = 50
$ws.Range("F47").Value = 25
$ws = $wb.Worksheets.Item("演出")
$ws.Range("F5").Value = 174
$ws.Range("F7").Value = 69
$ws.Range("F8").Value = 4
$ws.Range("F10").Value = 90
$ws.Range("F13").Value = 92
$ws.Range("G13").Value = 159
$ws.Range("F14").Value = 4396
$ws = $wb.Worksheets.Item("本地生活")
$ws.Range("F2").Value = 2136
$ws.Range("F3").Value = 591
$ws.Range("F4").Value = 510
$ws = $wb.Worksheets.Item("全部类型")
$ws.Range("F2").Value = 2136
$ws.Range("F3").Value = 591
$ws.Range("F4").Value = 10370
$ws.Range("F5").Value = 242
$ws.Range("F6").Value = 1141
$ws.Range("F7").Value = 510
$ws.Range("F8").Value = 1027
$ws.Range("F9").Value = 810
$ws.Range("F10").Value = 174
$ws.Range("F11").Value = 329
$ws.Range("F12").Value = 1078
$ws.Range("F14").Value = 143
$ws.Range("F15").Value = 833
$ws.Range("F16").Value = 362
$ws.Range("F17").Value = 1742
$ws.Range("F19").Value = 878
$ws.Range("F20").Value = 787
$ws.Range("F21").Value = 520
$ws.Range("F22").Value = 757
$ws.Range("F23").Value = 854
$ws.Range("F26").Value = 83
$ws.Range("F27").Value = 591
$ws.Range("F28").Value = 4
$ws.Range("F30").Value = 598
$ws.Range("F33").Value = 988
$ws.Range("F34").Value = 90
$ws.Range("F36").Value = 468
$ws.Range("F37").Value = 137
$ws.Range("F38").Value = 205
$ws.Range("F40").Value = 347
$ws.Range("F41").Value = 1358
$ws.Range("F42").Value = 395
$ws.Range("F43").Value = 112
$ws.Range("F45").Value = 41
$ws.Range("F46").Value = 69
$ws.Range("F47").Value = 61
